# Apply the changes described by the diff:
#  - Metadata!B2 (URL): corecl -> CoreCL in path segment
#  - Metadata!B8 (Date): updated timestamp
#  - Elements!Y7 (Binding Value Set): corecl -> CoreCL in path segment
#  - Elements sheet column widths tightened (cols 1-9,11,15,20-34,36)

$wb = $excel.ActiveWorkbook

# --- Metadata sheet text updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://hl7chile.cl/fhir/ig/CoreCL/StructureDefinition/RazonNOrealizarseInm"
$wsMeta.Range("B8").Value = "2022-12-12T20:08:16-03:00"

# --- Elements sheet text update ---
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("Y7").Value = "https://hl7chile.cl/fhir/ig/CoreCL/ValueSet/VSRazonNOTinm"

# --- Elements sheet column width updates ---
# (values chosen so the engine's stored width lands as close as possible
#  to the target width from the diff; widths for hidden columns also
#  re-assert Hidden = $true since writing ColumnWidth resets it)
$wsElem.Columns.Item(1).ColumnWidth = 18.166666666666668
$wsElem.Columns.Item(2).ColumnWidth = 10.333333333333334
$wsElem.Columns.Item(3).ColumnWidth = 6.833333333333333
$wsElem.Columns.Item(3).Hidden = $true
$wsElem.Columns.Item(4).ColumnWidth = 5.0
$wsElem.Columns.Item(4).Hidden = $true
$wsElem.Columns.Item(5).ColumnWidth = 3.8333333333333335
$wsElem.Columns.Item(6).ColumnWidth = 4.166666666666667
$wsElem.Columns.Item(7).ColumnWidth = 13.833333333333334
$wsElem.Columns.Item(8).ColumnWidth = 11.166666666666666
$wsElem.Columns.Item(9).ColumnWidth = 11.833333333333334
$wsElem.Columns.Item(11).ColumnWidth = 59.166666666666664
$wsElem.Columns.Item(15).ColumnWidth = 12.666666666666666
$wsElem.Columns.Item(20).ColumnWidth = 14.833333333333334
$wsElem.Columns.Item(21).ColumnWidth = 15.333333333333334
$wsElem.Columns.Item(22).ColumnWidth = 16.166666666666668
$wsElem.Columns.Item(23).ColumnWidth = 15.5
$wsElem.Columns.Item(24).ColumnWidth = 18.0
$wsElem.Columns.Item(25).ColumnWidth = 54.166666666666664
$wsElem.Columns.Item(26).ColumnWidth = 4.833333333333333
$wsElem.Columns.Item(27).ColumnWidth = 18.833333333333332
$wsElem.Columns.Item(28).ColumnWidth = 39.166666666666664
$wsElem.Columns.Item(29).ColumnWidth = 14.166666666666666
$wsElem.Columns.Item(30).ColumnWidth = 11.5
$wsElem.Columns.Item(31).ColumnWidth = 16.833333333333332
$wsElem.Columns.Item(31).Hidden = $true
$wsElem.Columns.Item(32).ColumnWidth = 8.666666666666666
$wsElem.Columns.Item(32).Hidden = $true
$wsElem.Columns.Item(33).ColumnWidth = 9.0
$wsElem.Columns.Item(33).Hidden = $true
$wsElem.Columns.Item(34).ColumnWidth = 11.333333333333334
$wsElem.Columns.Item(36).ColumnWidth = 21.833333333333332
